$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.275.44"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.21%  "

# Row 3 - Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.909.02"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.43%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.22%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.69"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.02%  "

# Row 6 - USDC
$ws.Range("E6").Value = "  +0.24%  "

# Row 7 - XRP
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5363"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +3.06%  "

# Row 8 - Cardano
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3814"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +1.13%  "

# Row 9 - Dogecoin
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07295"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.23%  "

# Row 10 - Solana
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.11"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +4.43%  "

# Row 11 - Polygon
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9034"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.03%  "

# Row 12 - TRON
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08206"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.88%  "

# Row 13 - Litecoin
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "95.86"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.61%  "

# Row 14 - Polkadot
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.349"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.33%  "

# Row 15 - BinanceUSD
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.002"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.13%  "

# Row 16 - Avalanche
$ws.Range("E16").Value = "  +1.78%  "

# Row 17 - ShibaInu
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008653"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.26%  "

# Row 18 - Dai
$ws.Range("E18").Value = "  +0.24%  "

# Row 19 - WrappedBTC
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "27.313.00"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.29%  "

# Rows 20/21 - Uniswap and WrappedEther swap positions
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.038"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.10%  "

$ws.Range("B21").Value = "WrappedEther"
$ws.Range("C21").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.107.60"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -41.76%  "

# Row 22 - Cosmos
$ws.Range("E22").Value = "  +1.38%  "

# Row 23 - Chainlink
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.522"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.52%  "

# Row 24 - Monero
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "149.79"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.87%  "

# Row 25 - LidoDAOToken
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.292"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.10%  "

# Row 26 - EthereumClassic
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "18.29"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.28%  "

# Row 27 - Toncoin
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.747"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.15%  "

# Row 28 - BitcoinCash
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "117.00"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.69%  "

# Row 29 - InternetComputer(DFINITY)
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.832"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.01%  "

# Row 30 - Filecoin
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.808"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.94%  "

# Row 31 - Stellar
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09281"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.36%  "

# Row 32 - ImmutableX
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.8336"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +4.24%  "

# Row 33 - Hedera
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05068"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.09%  "

# Row 34 - ARBITRUM
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.226"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -1.04%  "

# Row 35 - HuobiToken
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.998"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +1.79%  "

# Row 36 - MXToken
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.352"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -2.38%  "

# Row 37 - RenderToken
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.699"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +3.68%  "

# Row 38 - TheSandbox
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5776"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +1.05%  "

# Row 39 - VeChain
$ws.Range("E39").Value = "  +0.48%  "

# Row 40 - TrustWalletToken
$ws.Range("E40").Value = "  +0.10%  "

# Row 41 - Aptos
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "9.343"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +3.43%  "

# Row 42 - FraxShare
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.576"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.06%  "

# Row 43 - Quant
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "117.74"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +1.04%  "

# Row 44 - Algorand
$ws.Range("E44").Value = "  +0.48%  "

# Row 45 - Decentraland
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4933"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.55%  "

# Row 46 - PaxDollar
$ws.Range("E46").Value = "  +0.23%  "

# Row 47 - EnergySwap
$ws.Range("E47").Value = "  -0.02%  "

# Row 48 - NEARProtocol
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.641"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.87%  "

# Row 49 - Elrond
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "38.55"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +2.35%  "

# Row 50 - Cronos
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06114"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +2.68%  "

# Row 51 - Aave
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "63.27"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.95%  "
